# Add sourcefiles for OBJs 118, 119, 120, 121 (obj000117 - obj000119 rows)
# i.e. three new rows of data across the "Objects", "Images" and "Albums" sheets.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Objects": fill in rows 56-58 (cat-nos obj000117-01 .. obj000119-01
# already exist in column A) with the rest of their row data.
# ----------------------------------------------------------------------
$wsObjects = $wb.Worksheets.Item("Objects")

$objectsData = @(
    @{ Row = 56; B = "2015-08-09-obj000117-01.textile"; C = "art"; D = "artworks"; E = "Fabric/hanging"; F = "Untitled"; G = "pic000220"; H = "alb000117" },
    @{ Row = 57; B = "2015-08-09-obj000118-01.textile"; C = "art"; D = "artworks"; E = "Fabric/hanging"; F = "Untitled"; G = "pic000221"; H = "alb000118" },
    @{ Row = 58; B = "2015-08-09-obj000119-01.textile"; C = "art"; D = "artworks"; E = "Fabric/hanging"; F = "Untitled"; G = "pic000222"; H = "alb000119" }
)

foreach ($item in $objectsData) {
    $r = $item.Row
    $wsObjects.Cells.Item($r, 2).Value = $item.B
    $wsObjects.Cells.Item($r, 3).Value = $item.C
    $wsObjects.Cells.Item($r, 4).Value = $item.D
    $wsObjects.Cells.Item($r, 5).Value = $item.E
    $wsObjects.Cells.Item($r, 6).Value = $item.F
    $wsObjects.Cells.Item($r, 7).Value = $item.G
    $wsObjects.Cells.Item($r, 8).Value = $item.H
}

# ----------------------------------------------------------------------
# Sheet "Images": append rows 161-163 with the new picture records.
# Column C keeps the "photo number" formatting already used by the
# sibling cells above (copy that formatting across, same as the rest of
# the table), so we copy formats from the row right above the new block.
# ----------------------------------------------------------------------
$wsImages = $wb.Worksheets.Item("Images")

$imagesData = @(
    @{ Row = 161; A = "pic000220"; B = "2015-08-09-pic000220.textile"; C = "S-33-0021"; I = "alb000117" },
    @{ Row = 162; A = "pic000221"; B = "2015-08-09-pic000221.textile"; C = "S-32-0018"; I = "alb000118" },
    @{ Row = 163; A = "pic000222"; B = "2015-08-09-pic000222.textile"; C = "S-34-0065"; I = "alb000119" }
)

foreach ($item in $imagesData) {
    $r = $item.Row
    $wsImages.Cells.Item($r, 1).Value = $item.A
    $wsImages.Cells.Item($r, 2).Value = $item.B
    $wsImages.Cells.Item($r, 3).Value = $item.C
    $wsImages.Cells.Item($r, 9).Value = $item.I
}

# Re-use the "photo number" cell formatting (Arial 10) already applied to
# column C of row 160 for the three freshly added rows, matching the
# existing table styling instead of inventing new style entries.
$wsImages.Range("C160").Copy() | Out-Null
$wsImages.Range("C161:C163").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ----------------------------------------------------------------------
# Sheet "Albums": append the three new album ids.
# ----------------------------------------------------------------------
$wsAlbums = $wb.Worksheets.Item("Albums")
$wsAlbums.Cells.Item(24, 2).Value = "alb000117"
$wsAlbums.Cells.Item(25, 2).Value = "alb000118"
$wsAlbums.Cells.Item(26, 2).Value = "alb000119"

# ----------------------------------------------------------------------
# Update the active selections to match where the editor ended up after
# making the edits (bottom of each updated table).
# ----------------------------------------------------------------------
$wsImages.Activate()
$wsImages.Range("C163").Select() | Out-Null

$wsAlbums.Activate()
$wsAlbums.Range("B26").Select() | Out-Null

$wsObjects.Activate()
$wsObjects.Range("H58").Select() | Out-Null
